# Update Xbox GDK Samples to November GDK release.
$d = $word.ActiveDocument

function Replace-ExactRange($findText, $replaceText, $restyleHyperlink) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
        return
    }
    $start = $r.Start
    $r.Text = $replaceText
    if ($restyleHyperlink) {
        $newEnd = $start + $replaceText.Length
        $fix = $d.Range($start, $newEnd)
        $fix.Style = "Hyperlink"
    }
}

# 1. "If using Project Scarlett, set the active solution platform to " ->
#    3 runs: "If using " / "an Xbox Series X|S devkit" / ", set the active solution platform to "
Replace-ExactRange "Project Scarlett" "an Xbox Series X|S devkit" $false

# 2. First "SHGetKnownFolderPath" (split SHG/e/tKnownFolderPath) -> single run
Replace-ExactRange "SHGetKnownFolderPath" "SHGetKnownFolderPath" $true

# 3. ") - Win32 apps | Microsoft Docs" (split ...Micros/o/ft Docs) -> single run
Replace-ExactRange ") - Win32 apps | Micros" ") - Win32 apps | Microsoft Docs" $true

# 4. "GetTempPath" (split GetTe/m/pPath) -> single run
Replace-ExactRange "GetTempPath" "GetTempPath" $true

# 5. Second "SHGetKnownFolderPath" (split SHG/e/tKnownFolderPath) -> single run
Replace-ExactRange "SHGetKnownFolderPath" "SHGetKnownFolderPath" $true

# 6. "Knownfolders.h" (split Knownfolder/s/.h) -> single run
Replace-ExactRange "Knownfolders.h" "Knownfolders.h" $true

# 7. "GetDiskFreeSpa" (partial merge of GetDiskFreeSp/a only) -> single run
Replace-ExactRange "GetDiskFreeSp" "GetDiskFreeSpa" $true

# 8. "GetDiskFreeSpaceEx" (split GetDi/s/kFreeSpaceEx) -> single run
Replace-ExactRange "GetDiskFreeSpaceEx" "GetDiskFreeSpaceEx" $true
